# Append a new "filtered feed" row to the Filtered Feeds worksheet, mirroring
# what the upstream workflow does when it finds a new matching article.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = $ws.UsedRange.Rows.Count + 1

$link  = "https://www.360dx.com/cancer/acrivon-therapeutics-open-clia-certified-lab-run-diagnostics-targeted-therapies"
$kw    = "CDx"
$title = "Acrivon Therapeutics to Open CLIA-Certified Lab to Run Diagnostics for Targeted Therapies"

$aCell = $ws.Cells.Item($newRow, 1)
$bCell = $ws.Cells.Item($newRow, 2)
$cCell = $ws.Cells.Item($newRow, 3)

$aCell.Value = $link
$bCell.Value = $kw
$cCell.Value = $title

# Match the existing "link" column formatting (the Hyperlink cell style) used
# by every other row in column A.
$aCell.Style = $ws.Cells.Item($newRow - 1, 1).Style

# Turn the link cell into a real hyperlink, same as the other rows.
$ws.Hyperlinks.Add($aCell, $link) | Out-Null

# Re-apply the Hyperlink style in case adding the hyperlink changed it.
$aCell.Style = $ws.Cells.Item($newRow - 1, 1).Style
